$wb = $excel.ActiveWorkbook

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsProcess.Columns.Item(8).Delete()

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Columns.Item(8).Delete()
